$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the total "Valor Mora" amount
$ws.Range("E11").Value = 376144

# Update "Cant. Periodos" count (now only 1 period remains)
$ws.Range("F13").Value = 1

# Remove the two obsolete period rows (2504 and 2503), keeping only 2505
$ws.Rows("17:18").Delete()
